$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Extend header row 1 with two new bold/bordered header cells (P1=14, Q1=15)
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
# Copy the formatting (style) from O1, the last existing header cell, onto P1:Q1
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Swap the I/K and M/O columns for all data rows (2-25)
$ws.Range("I2:I25").Value = 2
$ws.Range("K2:K25").Value = 1
$ws.Range("M2:M25").Value = 2
$ws.Range("O2:O25").Value = 1

# Add new data columns P and Q (value 2) for all data rows (2-25)
$ws.Range("P2:Q25").Value = 2
